# Scheduled-runner update: refresh market-price snapshot values across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2922.7563
$ws.Range("I15").Value = 2922.7563
$ws.Range("K15").Value = 8768.268899999999
$ws.Range("M15").Value = -8599.268899999999
$ws.Range("H113").Value = 16668960
$ws.Range("J113").Value = 2776
$ws.Range("L113").Value = 2776
$ws.Range("N113").Value = -9284
$ws.Range("H129").Value = 676.617
$ws.Range("J129").Value = 863.76666
$ws.Range("L129").Value = 2591.29998
$ws.Range("N129").Value = -12591.29998
$ws.Range("H132").Value = 5957652
$ws.Range("I132").Value = 7756577
$ws.Range("J132").Value = 7361.846
$ws.Range("K132").Value = 23269731
$ws.Range("L132").Value = 22085.538
$ws.Range("M132").Value = -23267201
$ws.Range("N132").Value = -27145.538
$ws.Range("H137").Value = 1210.6825
$ws.Range("I137").Value = 878.1389
$ws.Range("J137").Value = 1654.0741
$ws.Range("K137").Value = 2634.4167
$ws.Range("L137").Value = 4962.2223
$ws.Range("M137").Value = -84.41670000000022
$ws.Range("N137").Value = -10062.2223
$ws.Range("H138").Value = 1210.01
$ws.Range("I138").Value = 567.6585
$ws.Range("J138").Value = 1656.3898
$ws.Range("K138").Value = 1702.9755
$ws.Range("L138").Value = 4969.1694
$ws.Range("M138").Value = 3437.0245
$ws.Range("N138").Value = -15249.1694

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7518.375
$ws.Range("I32").Value = 6717.119
$ws.Range("K32").Value = 6717.119
$ws.Range("M32").Value = -6430.119
$ws.Range("H61").Value = 32258838
$ws.Range("I61").Value = 41667390
$ws.Range("J61").Value = 942.8570999999999
$ws.Range("K61").Value = 41667390
$ws.Range("L61").Value = 942.8570999999999
$ws.Range("M61").Value = -41667178
$ws.Range("N61").Value = -1366.8571
$ws.Range("H97").Value = 328.3793
$ws.Range("I97").Value = 330.95834
$ws.Range("J97").Value = 316
$ws.Range("K97").Value = 330.95834
$ws.Range("L97").Value = 316
$ws.Range("M97").Value = 165.04166
$ws.Range("N97").Value = -1308
$ws.Range("H136").Value = 32258838
$ws.Range("I136").Value = 41667390
$ws.Range("J136").Value = 942.8570999999999
$ws.Range("K136").Value = 125002170
$ws.Range("L136").Value = 2828.5713
$ws.Range("M136").Value = -124999620
$ws.Range("N136").Value = -7928.5713
$ws.Range("H141").Value = 32086.334
$ws.Range("J141").Value = 32086.334
$ws.Range("L141").Value = 32086.334
$ws.Range("N141").Value = -42446.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 83334640
$ws.Range("I99").Value = 125001150
$ws.Range("J99").Value = 1625
$ws.Range("K99").Value = 125001150
$ws.Range("L99").Value = 1625
$ws.Range("M99").Value = -124999652
$ws.Range("N99").Value = -4621
$ws.Range("H107").Value = 942.0714
$ws.Range("I107").Value = 787.6286
$ws.Range("J107").Value = 1714.2858
$ws.Range("K107").Value = 787.6286
$ws.Range("L107").Value = 1714.2858
$ws.Range("M107").Value = 1132.3714
$ws.Range("N107").Value = -5554.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H31").Value = 2236.0833
$ws.Range("I31").Value = 2243.4
$ws.Range("K31").Value = 2243.4
$ws.Range("M31").Value = -1948.4
$ws.Range("H34").Value = 2236.0833
$ws.Range("I34").Value = 2243.4
$ws.Range("K34").Value = 2243.4
$ws.Range("M34").Value = -2041.4
$ws.Range("H125").Value = 14000
$ws.Range("J125").Value = 14000
$ws.Range("L125").Value = 14000
$ws.Range("N125").Value = -18920
$ws.Range("H132").Value = 4318
$ws.Range("I132").Value = 5327.2593
$ws.Range("J132").Value = 2371.5715
$ws.Range("K132").Value = 15981.7779
$ws.Range("L132").Value = 7114.7145
$ws.Range("M132").Value = -13451.7779
$ws.Range("N132").Value = -12174.7145
$ws.Range("H134").Value = 12821953
$ws.Range("I134").Value = 1497.8966
$ws.Range("J134").Value = 50001270
$ws.Range("K134").Value = 4493.6898
$ws.Range("L134").Value = 150003810
$ws.Range("M134").Value = -1958.6898
$ws.Range("N134").Value = -150008880
$ws.Range("H141").Value = 81042.734
$ws.Range("J141").Value = 81042.734
$ws.Range("L141").Value = 81042.734
$ws.Range("N141").Value = -91402.734

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 38461644
$ws.Range("I12").Value = 142857310
$ws.Range("J12").Value = 81.8421
$ws.Range("K12").Value = 428571930
$ws.Range("L12").Value = 245.5263
$ws.Range("M12").Value = -428571757
$ws.Range("N12").Value = -591.5263
$ws.Range("H131").Value = 20002994
$ws.Range("I131").Value = 90909610
$ws.Range("J131").Value = 3691.1025
$ws.Range("K131").Value = 272728830
$ws.Range("L131").Value = 11073.3075
$ws.Range("M131").Value = -272723790
$ws.Range("N131").Value = -21153.3075
$ws.Range("H133").Value = 3759.3
$ws.Range("I133").Value = 2006.6666
$ws.Range("J133").Value = 4068.5881
$ws.Range("K133").Value = 6019.9998
$ws.Range("L133").Value = 12205.7643
$ws.Range("M133").Value = -959.9997999999996
$ws.Range("N133").Value = -22325.7643
$ws.Range("H134").Value = 3296.0667
$ws.Range("I134").Value = 1671.0625
$ws.Range("J134").Value = 5153.2144
$ws.Range("K134").Value = 5013.1875
$ws.Range("L134").Value = 15459.6432
$ws.Range("M134").Value = 56.8125
$ws.Range("N134").Value = -25599.6432
$ws.Range("H136").Value = 2591.923
$ws.Range("I136").Value = 2004.2858
$ws.Range("J136").Value = 3277.5
$ws.Range("K136").Value = 6012.857400000001
$ws.Range("L136").Value = 9832.5
$ws.Range("M136").Value = -912.8574000000008
$ws.Range("N136").Value = -20032.5
$ws.Range("H137").Value = 24199912
$ws.Range("I137").Value = 68183630
$ws.Range("J137").Value = 8866.549999999999
$ws.Range("K137").Value = 204550890
$ws.Range("L137").Value = 26599.65
$ws.Range("M137").Value = -204545790
$ws.Range("N137").Value = -36799.64999999999
$ws.Range("H138").Value = 2781.1333
$ws.Range("I138").Value = 2764.625
$ws.Range("J138").Value = 2800
$ws.Range("K138").Value = 8293.875
$ws.Range("L138").Value = 8400
$ws.Range("M138").Value = -3153.875
$ws.Range("N138").Value = -18680
$ws.Range("H139").Value = 2174.1333
$ws.Range("I139").Value = 2462.15
$ws.Range("J139").Value = 1598.1
$ws.Range("K139").Value = 7386.450000000001
$ws.Range("L139").Value = 4794.299999999999
$ws.Range("M139").Value = -2246.450000000001
$ws.Range("N139").Value = -15074.3
$ws.Range("H140").Value = 3031.4363
$ws.Range("I140").Value = 3376
$ws.Range("J140").Value = 2863.8108
$ws.Range("K140").Value = 10128
$ws.Range("L140").Value = 8591.432400000002
$ws.Range("M140").Value = -4948
$ws.Range("N140").Value = -18951.4324
$ws.Range("H141").Value = 3009.2727
$ws.Range("I141").Value = 2016.9
$ws.Range("K141").Value = 6050.700000000001
$ws.Range("M141").Value = -870.7000000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9617553
$ws.Range("J122").Value = 41669132
$ws.Range("L122").Value = 125007396
$ws.Range("N122").Value = -125012296
$ws.Range("H132").Value = 2355.35
$ws.Range("I132").Value = 2016.069
$ws.Range("J132").Value = 3249.818
$ws.Range("K132").Value = 6048.207
$ws.Range("L132").Value = 9749.454000000002
$ws.Range("M132").Value = -3518.207
$ws.Range("N132").Value = -14809.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1617.5416
$ws.Range("I7").Value = 1571.5294
$ws.Range("J7").Value = 1729.2858
$ws.Range("K7").Value = 1571.5294
$ws.Range("L7").Value = 1729.2858
$ws.Range("M7").Value = -1459.5294
$ws.Range("N7").Value = -1953.2858
$ws.Range("H22").Value = 527.75
$ws.Range("I22").Value = 387.66666
$ws.Range("J22").Value = 667.8333
$ws.Range("K22").Value = 387.66666
$ws.Range("L22").Value = 667.8333
$ws.Range("M22").Value = -92.66665999999998
$ws.Range("N22").Value = -1257.8333
$ws.Range("H27").Value = 527.75
$ws.Range("I27").Value = 387.66666
$ws.Range("J27").Value = 667.8333
$ws.Range("K27").Value = 387.66666
$ws.Range("L27").Value = 667.8333
$ws.Range("M27").Value = -280.66666
$ws.Range("N27").Value = -881.8333
$ws.Range("H40").Value = 2084.742
$ws.Range("I40").Value = 1623.8334
$ws.Range("J40").Value = 3665
$ws.Range("K40").Value = 1623.8334
$ws.Range("L40").Value = 3665
$ws.Range("M40").Value = -1487.8334
$ws.Range("N40").Value = -3937
$ws.Range("H45").Value = 5499.5
$ws.Range("I45").Value = 5499.5
$ws.Range("K45").Value = 5499.5
$ws.Range("M45").Value = -5092.5
$ws.Range("H68").Value = 1953.0869
$ws.Range("I68").Value = 1945.6666
$ws.Range("J68").Value = 1979.8
$ws.Range("K68").Value = 1945.6666
$ws.Range("L68").Value = 1979.8
$ws.Range("M68").Value = -1196.6666
$ws.Range("N68").Value = -3477.8
$ws.Range("H71").Value = 1953.0869
$ws.Range("I71").Value = 1945.6666
$ws.Range("J71").Value = 1979.8
$ws.Range("K71").Value = 9728.333000000001
$ws.Range("L71").Value = 9899
$ws.Range("M71").Value = -5984.333000000001
$ws.Range("N71").Value = -17387
$ws.Range("H126").Value = 1617.5416
$ws.Range("I126").Value = 1571.5294
$ws.Range("J126").Value = 1729.2858
$ws.Range("K126").Value = 4714.5882
$ws.Range("L126").Value = 5187.857400000001
$ws.Range("M126").Value = -2244.5882
$ws.Range("N126").Value = -10127.8574
$ws.Range("H132").Value = 29907.473
$ws.Range("I132").Value = 1645.3684
$ws.Range("J132").Value = 61494.53
$ws.Range("K132").Value = 4936.1052
$ws.Range("L132").Value = 184483.59
$ws.Range("M132").Value = -2406.1052
$ws.Range("N132").Value = -189543.59
$ws.Range("H136").Value = 3570.2927
$ws.Range("I136").Value = 3916.1177
$ws.Range("J136").Value = 1890.5714
$ws.Range("K136").Value = 11748.3531
$ws.Range("L136").Value = 5671.7142
$ws.Range("M136").Value = -9198.3531
$ws.Range("N136").Value = -10771.7142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6100096.5
$ws.Range("I122").Value = 8067497.5
$ws.Range("J122").Value = 1153.3
$ws.Range("K122").Value = 24202492.5
$ws.Range("L122").Value = 3459.9
$ws.Range("M122").Value = -24200042.5
$ws.Range("N122").Value = -8359.9
$ws.Range("H126").Value = 50506788
$ws.Range("J126").Value = 1863.3334
$ws.Range("L126").Value = 5590.0002
$ws.Range("N126").Value = -10530.0002
$ws.Range("H136").Value = 759.9474
$ws.Range("I136").Value = 642.3043
$ws.Range("J136").Value = 940.3333
$ws.Range("K136").Value = 1926.9129
$ws.Range("L136").Value = 2820.9999
$ws.Range("M136").Value = 623.0871
$ws.Range("N136").Value = -7920.9999
